$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the original row 2 (dated 39400 / 2007), shifting all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

# After the shift, rows 2-5 no longer carry a y_1_forecast (column E) value, so clear any
# leftover values that moved into that range.
$ws.Range("E2:E5").ClearContents()

# Recomputed y_1_forecast values (column E) for the remaining forecast rows.
$ws.Range("E6").Value = 4.792854588620821
$ws.Range("E7").Value = 7.239454936865775
$ws.Range("E8").Value = 0.4361429468412448
$ws.Range("E9").Value = 1.392195163617171
$ws.Range("E10").Value = 1.493220091771108
$ws.Range("E11").Value = 2.622364272988187
$ws.Range("E12").Value = 3.113086948791377
$ws.Range("E13").Value = 2.757421718286168
$ws.Range("E14").Value = -1.373617952268746
$ws.Range("E15").Value = 1.324283050325015
$ws.Range("E16").Value = -0.5025420863900898
$ws.Range("E17").Value = -2.371854438773213
$ws.Range("E18").Value = -1.436963918858969
